# Apply the updated arithmetic equations to each cell of the single
# 20x5 table, in row-major order (matches how the document's cells are
# laid out / enumerated).
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
    "42+35="
    "73-29="
    "11+73="
    "50-27="
    "75-70="
    "63-56="
    "44+2="
    "56-19="
    "11+22="
    "26-4="
    "94-72="
    "71+17="
    "16+78="
    "21+36="
    "57+16="
    "61+20="
    "94-44="
    "45+16="
    "97+2="
    "20+69="
    "37+9="
    "86-84="
    "87-81="
    "43-39="
    "6+10="
    "85-46="
    "5+75="
    "83+14="
    "4+67="
    "28+11="
    "20+74="
    "20-15="
    "52-15="
    "23-2="
    "67+16="
    "83-1="
    "85-3="
    "93-85="
    "54-48="
    "93-27="
    "17+1="
    "57-50="
    "1+31="
    "40-31="
    "76-63="
    "63+22="
    "81+12="
    "55-16="
    "23+21="
    "95-54="
    "38+16="
    "31-12="
    "71-9="
    "79-58="
    "10+77="
    "48+0="
    "39-0="
    "58-4="
    "54-45="
    "25-4="
    "90-68="
    "65+26="
    "97-21="
    "81-63="
    "21+72="
    "24+66="
    "16+6="
    "54+42="
    "58-3="
    "37-37="
    "40-19="
    "42+22="
    "96-24="
    "60-0="
    "79-56="
    "89-39="
    "14+39="
    "92-54="
    "64-45="
    "53+18="
    "26+42="
    "62-14="
    "89-76="
    "51+15="
    "49+48="
    "20+19="
    "1+49="
    "97-78="
    "16+66="
    "30+22="
    "18+3="
    "59-27="
    "25+26="
    "35+63="
    "77+22="
    "24+61="
    "74-14="
    "81-67="
    "52-47="
    "9+51="
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count

if (($rows * $cols) -ne $newValues.Count) {
    throw "Table shape ($rows x $cols) does not match expected replacement count ($($newValues.Count))"
}

$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}

Write-Output ("Updated cells: " + $idx)
